$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 65
$ws.Range("A65").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B65").Value = "POPULATION, ESTIMATE_CLASSIFICATION, ESTIMATE_METHOD"
$ws.Range("C65").Value = "Consistency (C1)"
$ws.Range("D65").Value = "2024-12-02 13:55:13"
$ws.Range("E65").Value = 0.91
$ws.Range("F65").Value = 0.9946666666666667
$ws.Range("G65").Value = "EwertM"

# Row 66
$ws.Range("A66").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B66").Value = "AREA, ANALYSIS_YR, NATURAL_ADULT_SPAWNERS, NATURAL_JACK_SPAWNERS, NATURAL_SPAWNERS_TOTAL, ADULT_BROODSTOCK_REMOVALS, JACK_BROODSTOCK_REMOVALS, TOTAL_BROODSTOCK_REMOVALS, OTHER_REMOVALS, TOTAL_RETURN_TO_RIVER, EFFECTIVE_FEMALES, WEIGHTED_PCT_SPAWN, NO_INSPECTIONS_USED, ACT_ID, POP_ID, GFE_ID"
$ws.Range("C66").Value = "Accuracy (A1)"
$ws.Range("D66").Value = "2024-12-02 13:55:14"
$ws.Range("E66").Value = "no threshold"
$ws.Range("F66").Value = 1
$ws.Range("G66").Value = "EwertM"

# Row 67
$ws.Range("A67").Value = "Johnstone Strait and Strait of Georgia NuSEDS_20241004"
$ws.Range("B67").Value = "NATURAL_ADULT_SPAWNERS, NATURAL_JACK_SPAWNERS, NATURAL_SPAWNERS_TOTAL, ADULT_BROODSTOCK_REMOVALS, JACK_BROODSTOCK_REMOVALS, TOTAL_BROODSTOCK_REMOVALS, OTHER_REMOVALS, TOTAL_RETURN_TO_RIVER"
$ws.Range("C67").Value = "Accuracy (A2)"
$ws.Range("D67").Value = "2024-12-02 13:55:28"
$ws.Range("E67").Value = 1.5
$ws.Range("F67").Value = 1
$ws.Range("G67").Value = "EwertM"

# Row 68
$ws.Range("A68").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B68").Value = "All columns"
$ws.Range("C68").Value = "Accuracy (A3)"
$ws.Range("D68").Value = "2024-12-02 13:55:28"
$ws.Range("E68").Value = "no threshold"
$ws.Range("F68").Value = 1
$ws.Range("G68").Value = "EwertM"

# Row 69
$ws.Range("A69").Value = "Johnstone Strait and Strait of Georgia NuSEDS_20241004"
$ws.Range("B69").Value = "All columns"
$ws.Range("C69").Value = "Completeness (P)"
$ws.Range("D69").Value = "2024-12-02 13:55:41"
$ws.Range("E69").Value = 0.75
$ws.Range("F69").Value = 0.8482207305966877
$ws.Range("G69").Value = "EwertM"

# Row 70
$ws.Range("A70").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("B70").Value = "site_species_id, project_name, project_description, ecosystem_type, species_name, CU_Name, SMU_Display, SMU_ID"
$ws.Range("C70").Value = "Consistency (C1)"
$ws.Range("D70").Value = "2024-12-02 14:56:59"
$ws.Range("E70").Value = 0.91
$ws.Range("F70").Value = 0.9206204379562044
$ws.Range("G70").Value = "EwertM"

# Row 71
$ws.Range("A71").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B71").Value = "AREA, ANALYSIS_YR, NATURAL_ADULT_SPAWNERS, NATURAL_JACK_SPAWNERS, NATURAL_SPAWNERS_TOTAL, ADULT_BROODSTOCK_REMOVALS, JACK_BROODSTOCK_REMOVALS, TOTAL_BROODSTOCK_REMOVALS, OTHER_REMOVALS, TOTAL_RETURN_TO_RIVER, EFFECTIVE_FEMALES, WEIGHTED_PCT_SPAWN, NO_INSPECTIONS_USED, ACT_ID, POP_ID, GFE_ID"
$ws.Range("C71").Value = "Accuracy (A1)"
$ws.Range("D71").Value = "2024-12-02 14:57:00"
$ws.Range("E71").Value = "no threshold"
$ws.Range("F71").Value = 1
$ws.Range("G71").Value = "EwertM"

# Row 72
$ws.Range("A72").Value = "Johnstone Strait and Strait of Georgia NuSEDS_20241004"
$ws.Range("B72").Value = "NATURAL_ADULT_SPAWNERS, NATURAL_JACK_SPAWNERS, NATURAL_SPAWNERS_TOTAL, ADULT_BROODSTOCK_REMOVALS, JACK_BROODSTOCK_REMOVALS, TOTAL_BROODSTOCK_REMOVALS, OTHER_REMOVALS, TOTAL_RETURN_TO_RIVER"
$ws.Range("C72").Value = "Accuracy (A2)"
$ws.Range("D72").Value = "2024-12-02 14:57:11"
$ws.Range("E72").Value = 1.5
$ws.Range("F72").Value = 1
$ws.Range("G72").Value = "EwertM"

# Row 73
$ws.Range("A73").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B73").Value = "All columns"
$ws.Range("C73").Value = "Accuracy (A3)"
$ws.Range("D73").Value = "2024-12-02 14:57:12"
$ws.Range("E73").Value = "no threshold"
$ws.Range("F73").Value = 1
$ws.Range("G73").Value = "EwertM"

# Row 74
$ws.Range("A74").Value = "Johnstone Strait and Strait of Georgia NuSEDS_20241004"
$ws.Range("B74").Value = "All columns"
$ws.Range("C74").Value = "Completeness (P)"
$ws.Range("D74").Value = "2024-12-02 14:57:24"
$ws.Range("E74").Value = 0.75
$ws.Range("F74").Value = 0.8482207305966877
$ws.Range("G74").Value = "EwertM"

# Row 75
$ws.Range("A75").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("C75").Value = "Accuracy (A1)"
$ws.Range("D75").Value = "2024-12-02 14:59:48"
$ws.Range("E75").Value = "no threshold"
$ws.Range("G75").Value = "EwertM"

# Row 76
$ws.Range("A76").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("C76").Value = "Accuracy (A1)"
$ws.Range("D76").Value = "2024-12-02 15:00:24"
$ws.Range("E76").Value = "no threshold"
$ws.Range("G76").Value = "EwertM"

# Row 77
$ws.Range("A77").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("B77").Value = "site_latitude, site_longitude"
$ws.Range("C77").Value = "Accuracy (A2)"
$ws.Range("D77").Value = "2024-12-02 15:00:25"
$ws.Range("E77").Value = 1.5
$ws.Range("F77").Value = 1
$ws.Range("G77").Value = "EwertM"

# Row 78
$ws.Range("A78").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("B78").Value = "All columns"
$ws.Range("C78").Value = "Accuracy (A3)"
$ws.Range("D78").Value = "2024-12-02 15:00:25"
$ws.Range("E78").Value = "no threshold"
$ws.Range("F78").Value = 1
$ws.Range("G78").Value = "EwertM"

# Row 79
$ws.Range("A79").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("B79").Value = "All columns"
$ws.Range("C79").Value = "Completeness (P)"
$ws.Range("D79").Value = "2024-12-02 15:00:25"
$ws.Range("E79").Value = 0.75
$ws.Range("F79").Value = 0.7664092664092664
$ws.Range("G79").Value = "EwertM"

# Row 80
$ws.Range("A80").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("C80").Value = "Accuracy (A1)"
$ws.Range("D80").Value = "2024-12-02 15:01:13"
$ws.Range("E80").Value = "no threshold"
$ws.Range("G80").Value = "EwertM"

# Row 81
$ws.Range("A81").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("C81").Value = "Accuracy (A1)"
$ws.Range("D81").Value = "2024-12-02 15:01:47"
$ws.Range("E81").Value = "no threshold"
$ws.Range("G81").Value = "EwertM"

# Row 82
$ws.Range("A82").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("C82").Value = "Accuracy (A1)"
$ws.Range("D82").Value = "2024-12-02 15:02:09"
$ws.Range("E82").Value = "no threshold"
$ws.Range("G82").Value = "EwertM"

# Row 83
$ws.Range("A83").Value = "SalmonHabitatRestorationProjects_DataPortal_June_FinalFields_20240613"
$ws.Range("C83").Value = "Accuracy (A1)"
$ws.Range("D83").Value = "2024-12-02 15:02:33"
$ws.Range("E83").Value = "no threshold"
$ws.Range("G83").Value = "EwertM"

# Row 84
$ws.Range("A84").Value = "Pacific-Recreational-Fishery-Salmon-Head-Depots"
$ws.Range("B84").Value = "DEPOT NAME / NOM DU DÉPÔT, AREA / LA RÉGION, MUNICIPALITY / MUNICIPALITÉ, ADDRESS / ADRESSE, STORAGE INFORMATION / DÉTAILS DE STOCKAGE"
$ws.Range("C84").Value = "Consistency (C1)"
$ws.Range("D84").Value = "2024-12-02 15:31:28"
$ws.Range("E84").Value = 0.91
$ws.Range("F84").Value = 1
$ws.Range("G84").Value = "EwertM"

# Row 85
$ws.Range("A85").Value = "Pacific-Recreational-Fishery-Salmon-Head-Depots"
$ws.Range("B85").Value = "LATITUDE / LATITUDE, LONGITUDE / LONGITUDE"
$ws.Range("C85").Value = "Accuracy (A1)"
$ws.Range("D85").Value = "2024-12-02 15:31:28"
$ws.Range("E85").Value = "no threshold"
$ws.Range("F85").Value = 0.5
$ws.Range("G85").Value = "EwertM"

# Row 86
$ws.Range("A86").Value = "Pacific-Recreational-Fishery-Salmon-Head-Depots"
$ws.Range("B86").Value = "LATITUDE / LATITUDE, LONGITUDE / LONGITUDE"
$ws.Range("C86").Value = "Accuracy (A2)"
$ws.Range("D86").Value = "2024-12-02 15:31:28"
$ws.Range("E86").Value = 1.5
$ws.Range("F86").Value = 1
$ws.Range("G86").Value = "EwertM"

# Row 87
$ws.Range("A87").Value = "Pacific-Recreational-Fishery-Salmon-Head-Depots"
$ws.Range("B87").Value = "All columns"
$ws.Range("C87").Value = "Accuracy (A3)"
$ws.Range("D87").Value = "2024-12-02 15:31:28"
$ws.Range("E87").Value = "no threshold"
$ws.Range("F87").Value = 1
$ws.Range("G87").Value = "EwertM"

# Row 88
$ws.Range("A88").Value = "Pacific-Recreational-Fishery-Salmon-Head-Depots"
$ws.Range("B88").Value = "All columns"
$ws.Range("C88").Value = "Completeness (P)"
$ws.Range("D88").Value = "2024-12-02 15:31:28"
$ws.Range("E88").Value = 0.75
$ws.Range("F88").Value = 0.9764764764764765
$ws.Range("G88").Value = "EwertM"

